$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.348.12'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.62%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.935.63'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.01%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '252.18'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.94%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.7257'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.75%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.000'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3305'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.71%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '27.94'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +7.15%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07245'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.67%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8090'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.61%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08106'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.21%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.936.08'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.99%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.476'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.43%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '94.82'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.31%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.11'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.47%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.341.02'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.56%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008239'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +5.76%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '253.20'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.33%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.827'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.71%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.189.53'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.69%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.23%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.19%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.960'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.56%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.761'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.16%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '166.13'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.84%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.348'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +6.64%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.32'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.81%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1300'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.69%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.354'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.71%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.547'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.26%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.437'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.27%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.215'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.39%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05247'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.62%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.268'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +7.28%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7512'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.47%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.772'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.95%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01972'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.13%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.803'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.76%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '79.45'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.19%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.450'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.43%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4549'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.09%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.032'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.26%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8445'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.48%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.001'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.19%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '102.06'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.46%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.824'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.04%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.446'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.26%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '36.79'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.26%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4201'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.15%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06049'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.19%'
